# Add "Profit_gen" and "Profit_wind" worksheets at the end of the workbook,
# matching the commit "Profit of generators and wind in step 3".

$wb = $excel.ActiveWorkbook

# --- Profit_gen sheet (appended after the last existing sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGen = $wb.Worksheets.Add($null, $lastSheet)
$wsGen.Name = "Profit_gen"

$genHeaders = @("x1","x2","x3","x4","x5","x6","x7","x8","x9","x10","x11","x12")
for ($c = 0; $c -lt $genHeaders.Length; $c++) {
    $wsGen.Cells.Item(1, $c + 1).Value = $genHeaders[$c]
}

$genData = @(
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,316.64547677261635,0.0,107.99823960880212,-0.8647432762836615),
    @(0.0,0.0,0.0,0.0,0.0,28.030733679950345,0.0,1592.0,1751.318313113808,0.0,0.0,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,788.0,813.0847457627124,0.0,0.0,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1785.62095531587,0.0,114.70000000000027,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1771.9383667180268,0.0,0.0,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,106.22739018087941,0.0,0.0,0.0,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1771.9383667180277,0.0,0.0,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1785.62095531587,0.0,114.70000000000027,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1785.620955315871,0.0,114.70000000000027,0.0),
    @(0.0,0.0,0.0,0.0,0.0,55.482622838555244,0.0,1592.0,1742.6666666666665,0.0,114.70000000000027,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1759.1217967590987,0.0,10.907739788405252,-15.565331459197523),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1869.1182795698933,0.0,218.49226021159575,95.92062968752907),
    @(272.30120547945353,371.1104986301377,0.0,0.0,0.0,406.73656164383374,0.0,222.6712328767112,0.0,0.0,114.70000000000027,0.0),
    @(0.0,-74.03806594040384,0.0,0.0,0.0,59.86282990475593,0.0,1592.0,1742.6666666666665,0.0,123.77189552377104,10.242462688127944),
    @(0.0,0.0,0.0,0.0,0.0,55.482622838555244,0.0,1592.0,1742.6666666666665,0.0,112.94305027677683,-0.6234337727566981),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,79.87074829931998,0.0,0.0,117.16817359243396,2.276778293979987),
    @(0.0,0.0,0.0,0.0,0.0,77.94232397148016,20.213731019632633,1592.0,1725.2785109508532,0.0,112.23182640756659,-0.3659342224186162),
    @(0.0,0.0,0.0,0.0,0.0,29.620238826888453,0.0,1592.0,1750.4304073585681,0.0,2.3587084520722783,-30.84907953626623),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1858.8938271604939,0.0,227.04129154792827,117.82298267014585),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1820.0,0.0,114.70000000000027,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1802.4209553158707,0.0,114.70000000000027,0.0),
    @(0.0,0.0,0.0,0.0,0.0,58.39171171171142,0.0,1592.0,1742.6666666666665,0.0,114.70000000000027,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1785.6209553158715,0.0,114.70000000000027,0.0),
    @(0.0,0.0,0.0,0.0,0.0,0.0,0.0,1592.0,1766.9764252035998,0.0,0.0,0.0)
)
for ($r = 0; $r -lt $genData.Length; $r++) {
    $row = $genData[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsGen.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# --- Profit_wind sheet (appended after Profit_gen) ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWind = $wb.Worksheets.Add($null, $lastSheet2)
$wsWind.Name = "Profit_wind"

$windHeaders = @("x1","x2","x3","x4","x5","x6")
for ($c = 0; $c -lt $windHeaders.Length; $c++) {
    $wsWind.Cells.Item(1, $c + 1).Value = $windHeaders[$c]
}

$windData = @(
    @(0.0,1764.1645106693788,24.56961514124597,59.71755929916871,127.08651232915894,66.59913855616138),
    @(0.0,0.0,57.06788892703634,106.29970933910111,270.8297164496531,166.36487431802777),
    @(0.0,273.71214120050894,67.52417745099999,101.1217010128746,303.11297057659317,150.6157284648102),
    @(603.7008952602121,565.4896077852261,81.3636065392419,132.78213915675806,364.4929144614483,229.82430895691834),
    @(35.152804974970806,721.4715597694662,82.88974384899849,163.64768352234205,370.8949244992298,239.8082906082434),
    @(553.1258632030689,509.6813907170629,68.55834072377009,88.5196783832,315.7506225320713,127.3598853348),
    @(736.7149796399449,850.8752272215812,72.90980180423728,158.32526935747305,393.63137894375956,228.09853001232665),
    @(2324.686195821684,2525.5438299487514,79.18315579323271,169.67829055176574,404.45377960018504,244.8346014185639),
    @(2261.688907290222,2369.046941785402,88.02688520599075,184.0959060237658,384.7155687432974,248.18868668679198),
    @(2188.8049364883227,1690.5447426749615,87.6507466104,174.92877461546666,393.617401982,221.5342059084),
    @(1027.9318716596445,507.42856893374227,80.51940419542758,184.78093259913416,425.4533482773728,218.18539890617083),
    @(3068.7019946498112,880.9958888170329,61.38658736221934,166.8825564239656,369.83850717187084,232.74593811375487),
    @(3535.4837346624263,670.2285187503554,62.00012938403007,83.78086662999999,354.6703556478489,122.00111991809999),
    @(2413.1983214506836,536.6913256903616,72.3148889124,142.5762548144,393.19459055799996,189.6923788584),
    @(1697.5274818980781,924.4218857102076,71.699423992,139.8161149352,391.75955052,143.5620925464),
    @(2581.2339400956103,2042.212990968979,55.55708445063268,79.5734814586,302.3876695817348,89.45287411289999),
    @(1207.8420700213378,1729.3274248472937,69.5631720453778,143.341538527427,371.96259277632697,141.56280425045657),
    @(497.34267158347666,716.2908989043865,70.43685818952872,145.13721899865888,373.85750972217227,132.96338508450998),
    @(2785.5432192414805,3122.1419819364464,47.35125585571851,158.50941891914272,320.1337230693778,168.05367503624888),
    @(2536.89164316,2948.0605530884995,60.05645001720001,153.9075104196,368.43240295399994,197.1204255432),
    @(999.2420089324289,2695.88155307651,69.81371788497442,146.1744108914071,411.3536329917461,203.44053090802623),
    @(2059.5467275024857,939.0368205995104,85.8990962316,152.26321811946667,399.54399073599996,188.83870320920002),
    @(2025.5156256020166,2013.2748130493144,83.25148657916792,150.6486088693005,409.0565846319876,95.22081841556857),
    @(592.6046382804761,824.3680575480814,83.51396095058809,127.20775752661977,408.45884784385777,47.524835987033)
)
for ($r = 0; $r -lt $windData.Length; $r++) {
    $row = $windData[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsWind.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
